# Celerio-style document.xlsx update:
#   - On the "List" sheet, the documentBinary / documentFileName columns move
#     to sit right after accountId (before documentContentType), instead of
#     after documentSize.
#   - On the "Search" sheet, the trailing "document_account" / "account"
#     placeholders are replaced with "search_full_text" / "search_full_text".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "List": re-order the label row (row 1) and the value row (row 2).
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("List")

# New header order for columns C..F
$wsList.Range("C1").Value = "`${msg.getProperty('document_documentBinary')}"
$wsList.Range("D1").Value = "`${msg.getProperty('document_documentFileName')}"
$wsList.Range("E1").Value = "`${msg.getProperty('document_documentContentType')}"
$wsList.Range("F1").Value = "`${msg.getProperty('document_documentSize')}"

# New value order for columns C..F
$wsList.Range("C2").Value = "`${document.documentBinary}"
$wsList.Range("D2").Value = "`${document.documentFileName}"
$wsList.Range("E2").Value = "`${document.documentContentType}"
$wsList.Range("F2").Value = "`${document.documentSize}"

# ---------------------------------------------------------------------------
# Sheet "Search": replace the document_account / account row with
# search_full_text / search_full_text.
# ---------------------------------------------------------------------------
$wsSearch = $wb.Worksheets.Item("Search")

$wsSearch.Range("A4").Value = "`${msg.getProperty('search_full_text')}"
$wsSearch.Range("B4").Value = "`${search_full_text}"
